$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.641.36"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.791.20"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4939"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3828"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09251"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.086"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.207"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "1.790.90"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.100"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001099"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06532"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.886"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "27.673.96"
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.222"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.15%  "
$ws.Range("D27").Value = "2.001.07"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.369"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1063"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.040"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.605"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.480"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.00%  "
$ws.Range("E35").Value = "  -6.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.748"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02278"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2108"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.871"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6080"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.140"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5786"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.647"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.18%  "
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.268"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.909"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.161"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06698"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "
